# Edit the FlashScore games workbook:
#  - Header row: the "Odd_CS_4-4" label moves from AM1 to AG1, and the labels
#    that used to sit in AG1:AL1 (0-1,0-2,1-2,0-3,1-3,2-3) each shift one
#    column to the right, into AH1:AM1. Everything from AN1 onward is
#    untouched.
#  - The first data row (Colombia - Bucaramanga vs Fortaleza) is removed.
#  - The remaining data row (USA - New Mexico vs Las Vegas Lights) becomes
#    row 2, with refreshed odds values (including its "Odd_CS_4-4" value
#    now living in the new AG position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix up the header row labels for columns AG..AM.
$ws.Range("AG1").Value = "Odd_CS_4-4"
$ws.Range("AH1").Value = "Odd_CS_0-1"
$ws.Range("AI1").Value = "Odd_CS_0-2"
$ws.Range("AJ1").Value = "Odd_CS_1-2"
$ws.Range("AK1").Value = "Odd_CS_0-3"
$ws.Range("AL1").Value = "Odd_CS_1-3"
$ws.Range("AM1").Value = "Odd_CS_2-3"

# 2) Drop the old first data row; the remaining data row moves up to row 2.
$ws.Rows(2).Delete(-4162)  # xlShiftUp

# 3) Refresh row 2 (the USA match) with its updated values. The date (B2)
#    already reads "09/11/2024" for both the removed and the surviving row,
#    so it is left untouched (re-assigning it would coerce the text into a
#    real date serial, which the source file never uses).
$ws.Range("A2").Value = "tSZqiGYq"
$ws.Range("C2").Value = "23:30"
$ws.Range("D2").Value = "USA - USL CHAMPIONSHIP"
$ws.Range("E2").Value = "New Mexico"
$ws.Range("F2").Value = "Las Vegas Lights"

$values = [ordered]@{
    "G2" = 2.12; "H2" = 3.5; "I2" = 2.95; "J2" = 2.67; "K2" = 2.25; "L2" = 3.45;
    "M2" = 1.04; "N2" = 8.25; "O2" = 1.24; "P2" = 3.7; "Q2" = 1.72; "R2" = 2.05;
    "S2" = 1.33; "T2" = 3.05; "U2" = 1.62; "V2" = 2.18; "W2" = 9; "X2" = 11.25;
    "Y2" = 8.75; "Z2" = 20; "AA2" = 16; "AB2" = 23; "AC2" = 8.25; "AD2" = 7;
    "AE2" = 12.5; "AF2" = 50; "AG2" = 300; "AH2" = 11.25; "AI2" = 17; "AJ2" = 10.75;
    "AK2" = 37; "AL2" = 23; "AM2" = 28; "AN2" = 4.25; "AO2" = 10.75; "AP2" = 17.5;
    "AQ2" = 40; "AR2" = 65; "AS2" = 200; "AT2" = 3.05; "AU2" = 6.7; "AV2" = 50;
    "AW2" = 5.1; "AX2" = 15.5; "AY2" = 20; "AZ2" = 70; "BA2" = 90; "BB2" = 200;
    "BC2" = 51; "BD2" = 51
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
